# TradingModel - 2021/11/15 data update
# Appends the day's trading-history rows (2021/11/15) to the "交易記錄" sheet
# and leaves the view scrolled/selected where the author left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date serial 44515 == 2021-11-15. Reuse the date cell format already used
# lower in the table (col A, rows 8-20 use the custom m"月"d"日" format).
$dateFormat = $ws.Cells.Item(20, 1).NumberFormat

# Columns: A Date | B Stock_Id | C Action | D PositionSize | E Price
$newRows = @(
    @(44515, 1711, "short", -220, 28),
    @(44515, 3122, "long",    75,  65),
    @(44515, 3033, "short", -180, 31.15),
    @(44515, 6138, "long",    30, 203),
    @(44515, 2436, "long",    10, 105.5),
    @(44515, 3122, "long",    15,  66),
    @(44515, 3189, "short",  -27, 246),
    @(44515, 6271, "long",    20, 302.5)
)

$row = 21
foreach ($r in $newRows) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 1).NumberFormat = $dateFormat
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $row = $row + 1
}

# Match the author's final view/selection state.
[void]$ws.Range("A12").Select()
[void]$ws.Range("B23:E23").Select()
